# NIT-9015420606.xlsx — "Elimina EC anteriores y se agregan nuevos, se
# modifica base de datos"
#
# The mora (overdue-payment) detail table in B16:J36 is re-sorted: the
# worker records are regrouped by worker (CINTHIA LORETTE GONZALEZ PEREZ
# first, then ANGIE OCHOA NAVAS) with each worker's periods listed
# newest-to-oldest. The underlying set of (doc type, doc number, name,
# period, valor mora, salario basico) tuples is unchanged — this is a
# pure re-sort/re-grouping of the existing 21 detail rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=16; TipoDoc="CC"; NumDoc="1002320148"; Nombre="CINTHIA LORETTE GONZALEZ PEREZ"; Periodo="2406"; ValorMora=19488;  Salario=1218000 },
    @{ Row=17; TipoDoc="CC"; NumDoc="1002320148"; Nombre="CINTHIA LORETTE GONZALEZ PEREZ"; Periodo="2405"; ValorMora=48720;  Salario=1218000 },
    @{ Row=18; TipoDoc="CC"; NumDoc="1002320148"; Nombre="CINTHIA LORETTE GONZALEZ PEREZ"; Periodo="2404"; ValorMora=48720;  Salario=1218000 },
    @{ Row=19; TipoDoc="CC"; NumDoc="1002320148"; Nombre="CINTHIA LORETTE GONZALEZ PEREZ"; Periodo="2403"; ValorMora=48720;  Salario=1218000 },
    @{ Row=20; TipoDoc="CC"; NumDoc="1002320148"; Nombre="CINTHIA LORETTE GONZALEZ PEREZ"; Periodo="2402"; ValorMora=48720;  Salario=1218000 },
    @{ Row=21; TipoDoc="CC"; NumDoc="1002320148"; Nombre="CINTHIA LORETTE GONZALEZ PEREZ"; Periodo="2401"; ValorMora=48720;  Salario=1218000 },
    @{ Row=22; TipoDoc="CC"; NumDoc="1002320148"; Nombre="CINTHIA LORETTE GONZALEZ PEREZ"; Periodo="2312"; ValorMora=48720;  Salario=1218000 },
    @{ Row=23; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2406"; ValorMora=21101;  Salario=1300000 },
    @{ Row=24; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2405"; ValorMora=52753;  Salario=1300000 },
    @{ Row=25; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2404"; ValorMora=52753;  Salario=1300000 },
    @{ Row=26; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2403"; ValorMora=52753;  Salario=1300000 },
    @{ Row=27; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2402"; ValorMora=52753;  Salario=1300000 },
    @{ Row=28; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2401"; ValorMora=52753;  Salario=1300000 },
    @{ Row=29; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2312"; ValorMora=52753;  Salario=1300000 },
    @{ Row=30; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2311"; ValorMora=52753;  Salario=1300000 },
    @{ Row=31; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2310"; ValorMora=52753;  Salario=1300000 },
    @{ Row=32; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2309"; ValorMora=52753;  Salario=1300000 },
    @{ Row=33; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2308"; ValorMora=52753;  Salario=1300000 },
    @{ Row=34; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2307"; ValorMora=52753;  Salario=1300000 },
    @{ Row=35; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2306"; ValorMora=52753;  Salario=1300000 },
    @{ Row=36; TipoDoc="CC"; NumDoc="1128054473"; Nombre="ANGIE OCHOA NAVAS";              Periodo="2305"; ValorMora=52753;  Salario=1300000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 2).Value = $r.TipoDoc     # B - Tipo Doc Trabajador
    $ws.Cells.Item($n, 3).Value = $r.NumDoc       # C - N Doc Trabajador
    $ws.Cells.Item($n, 4).Value = $r.Nombre       # D - Nombre Trabajador
    $ws.Cells.Item($n, 5).Value = $r.Periodo      # E - Periodo Mora
    $ws.Cells.Item($n, 6).Value = $r.ValorMora    # F - Valor Mora
    $ws.Cells.Item($n, 7).Value = $r.Salario      # G - Salario Basico
}
